$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BME_RENAL_SPH")

# Delete the first four data rows (fiscal years 1313, 1413, 1513, 1613)
$ws.Rows("2:5").Delete()

# The remaining 4 rows (now rows 2-5) drop the E/F (resp. J/K) terms from
# their "actual" (resp. "budgeted") partial-OH formulas.
$ws.Range("B2:B5").Formula = "=C2-D2"
$ws.Range("G2:G5").Formula = "=H2-I2"
